$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.206.30"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "1.902.14"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5208"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07275"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9048"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08250"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.82%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.949.82"
$ws.Range("E13").Value = "  +3.40%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.289"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008671"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "27.244.90"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.093"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.03%  "
$ws.Range("D22").Value = "2.161.68"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.435"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.326"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.748"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.837"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.902"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09256"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05079"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.8001"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.957"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.589"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5711"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.080"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.025"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.591"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "116.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4860"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.626"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
